{"js": "// Load all paragraphs in the body so we can locate our anchor points by text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the anchor paragraphs we need to insert around / edit.\nlet socialMedias = null; // \"Social medias\" -> insert \"Or external paypal\" after this (same level)\nlet donateItem = null;   // \"Donate\" -> sanity check, new para goes right before this\nlet tbdItem = null;      // \"TBD: maybe google calendar implementation\" -> replace text + append sub items\n\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t === \"Social medias\") {\n    socialMedias = p;\n  } else if (t === \"Donate\") {\n    donateItem = p;\n  } else if (t.indexOf(\"TBD: maybe google calendar implementation\") !== -1) {\n    tbdItem = p;\n  }\n}\n\nif (!socialMedias || !donateItem || !tbdItem) {\n  throw new Error(\"Could not locate expected anchor paragraphs in the document.\");\n}\n\n// 1) Insert a new second-level bullet \"Or external paypal\" right after\n//    \"Social medias\" (i.e. right before \"Donate\"), matching its list level.\nconst payPalPara = socialMedias.insertParagraph(\"Or external \", Word.InsertLocation.after);\npayPalPara.insertText(\"paypal\", Word.InsertLocation.end);\n\n// 2) Replace the \"TBD: maybe google calendar implementation\" text with the\n//    new copy.\ntbdItem.insertText(\"Events that alumni would be interested in\", Word.InsertLocation.replace);\n\n// 3) Append the new sub-bullets about events, at increasing/decreasing list\n//    levels, after the paragraph we just edited.\nlet cursor = tbdItem;\n\nconst newItems = [\n  { text: \"Derby Days\", level: 2 },\n  { text: \"Black and White\", level: 3 },\n  { text: \"Sweetheart\", level: 2 },\n  { text: \"Sigtoberfest\", level: 2 },\n  { text: \"Email sign up for alumni\", level: 0 },\n  { text: \"Member page \\u2013 maybe just executive board\", level: 0 },\n  { text: \"Feed of Ig posts\", level: 0 },\n];\n\nfor (const item of newItems) {\n  cursor = cursor.insertParagraph(item.text, Word.InsertLocation.after);\n  cursor.listItem.level = item.level;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-ParaText($p) {\n    return $p.Range.Text.TrimEnd([char]13, [char]10, [char]7)\n}\n\nfunction Find-ParaIndex($d, $text, $exact) {\n    $i = 0\n    $found = -1\n    foreach ($p in $d.Paragraphs) {\n        $i = $i + 1\n        $t = (Get-ParaText $p).Trim()\n        if ($exact) {\n            if ($t -eq $text) { $found = $i }\n        } else {\n            if ($t.Contains($text)) { $found = $i }\n        }\n    }\n    return $found\n}\n\n# --- 1) Insert \"Or external paypal\" right before the \"Donate\" bullet,\n#        i.e. right after \"Social medias\" (same list level). ---\n$socialIdx = Find-ParaIndex $d \"Social medias\" $true\nif ($socialIdx -eq -1) { throw \"Could not locate the 'Social medias' paragraph.\" }\n$socialPara = $d.Paragraphs.Item($socialIdx)\n$socialPara.Range.InsertParagraphAfter()\n$payPalPara = $d.Paragraphs.Item($socialIdx + 1)\n$payPalPara.Range.Text = \"Or external paypal\"\n\n# --- 2) Replace \"TBD: maybe google calendar implementation\" text. ---\n$tbdIdx = Find-ParaIndex $d \"TBD: maybe google calendar implementation\" $false\nif ($tbdIdx -eq -1) { throw \"Could not locate the 'TBD: maybe google calendar implementation' paragraph.\" }\n$tbdPara = $d.Paragraphs.Item($tbdIdx)\n$tbdPara.Range.Text = \"Events that alumni would be interested in\"\n\n# --- 3) Append the new sub-bullets about events after that paragraph. ---\n$newItems = @(\n    @{ Text = \"Derby Days\"; Level = 3 },\n    @{ Text = \"Black and White\"; Level = 4 },\n    @{ Text = \"Sweetheart\"; Level = 3 },\n    @{ Text = \"Sigtoberfest\"; Level = 3 },\n    @{ Text = \"Email sign up for alumni\"; Level = 1 },\n    @{ Text = \"Member page \" + [char]0x2013 + \" maybe just executive board\"; Level = 1 },\n    @{ Text = \"Feed of Ig posts\"; Level = 1 }\n)\n\n$cursorIdx = $tbdIdx\nforeach ($item in $newItems) {\n    $cursorPara = $d.Paragraphs.Item($cursorIdx)\n    $cursorPara.Range.InsertParagraphAfter()\n    $cursorIdx = $cursorIdx + 1\n    $newPara = $d.Paragraphs.Item($cursorIdx)\n    $newPara.Range.Text = $item.Text\n    $newPara.Range.ListFormat.ListLevelNumber = $item.Level\n}\n"}
